# Removing remaining knock outs from iFerment
# Missed three additional reactions that should have been commented in
# iFerment, since knock outs are handled in assignments code.
# Set the corresponding B-column reaction flux values to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14:B19").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B25:B33").Value = 0
$ws.Range("B40").Value = 0
$ws.Range("B42").Value = 0
$ws.Range("B44:B45").Value = 0
$ws.Range("B77").Value = 0
$ws.Range("B119:B122").Value = 0
$ws.Range("B127").Value = 0
$ws.Range("B129:B130").Value = 0
$ws.Range("B136:B141").Value = 0
$ws.Range("B143:B145").Value = 0
$ws.Range("B147:B148").Value = 0
$ws.Range("B151").Value = 0
$ws.Range("B153").Value = 0
$ws.Range("B170").Value = 0
$ws.Range("B173").Value = 0
$ws.Range("B175").Value = 0
$ws.Range("B188").Value = 0
$ws.Range("B192").Value = 0
$ws.Range("B196").Value = 0
$ws.Range("B198").Value = 0
$ws.Range("B200:B201").Value = 0
$ws.Range("B210").Value = 0
$ws.Range("B226:B229").Value = 0
